$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "33.368.82"
$ws.Range("E2").Value = "  +11.57%  "
$ws.Range("D3").Value = "1.729.49"
$ws.Range("E3").Value = "  +4.76%  "
$ws.Range("D4").Value = "0.982"
$ws.Range("E4").Value = "  -1.42%  "
$ws.Range("D5").Value = "224.44"
$ws.Range("E5").Value = "  +3.38%  "
$ws.Range("E6").Value = "  +4.32%  "
$ws.Range("D7").Value = "0.993"
$ws.Range("E7").Value = "  -0.50%  "
$ws.Range("D8").Value = "31.01"
$ws.Range("E8").Value = "  +6.48%  "
$ws.Range("D9").Value = "44.88"
$ws.Range("E9").Value = "  -0.34%  "
$ws.Range("D10").Value = "0.272"
$ws.Range("E10").Value = "  +4.14%  "
$ws.Range("D11").Value = "0.0657"
$ws.Range("E11").Value = "  +7.55%  "
$ws.Range("D12").Value = "0.0914"
$ws.Range("E12").Value = "  +1.54%  "
$ws.Range("D13").Value = "1.973.17"
$ws.Range("E13").Value = "  +5.18%  "
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.714.04"
$ws.Range("E14").Value = "  +4.09%  "
$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").Value = "10.45"
$ws.Range("E15").Value = "  +5.69%  "
$ws.Range("D16").Value = "0.622"
$ws.Range("E16").Value = "  +4.27%  "
$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "32.952.80"
$ws.Range("E17").Value = "  +10.16%  "
$ws.Range("B18").Value = "Polkadot"
$ws.Range("C18").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D18").Value = "4.21"
$ws.Range("E18").Value = "  +7.32%  "
$ws.Range("D19").Value = "67.71"
$ws.Range("E19").Value = "  +4.98%  "
$ws.Range("D20").Value = "255.27"
$ws.Range("E20").Value = "  +6.58%  "
$ws.Range("D21").Value = "0.0₃0733"
$ws.Range("E21").Value = "  +3.65%  "
$ws.Range("D22").Value = "0.981"
$ws.Range("E22").Value = "  -1.66%  "
$ws.Range("D23").Value = "10.27"
$ws.Range("E23").Value = "  +3.05%  "
$ws.Range("D24").Value = "4.29"
$ws.Range("E24").Value = "  +3.19%  "
$ws.Range("D25").Value = "2.14"
$ws.Range("E25").Value = "  -1.62%  "
$ws.Range("D26").Value = "157.90"
$ws.Range("E26").Value = "  +0.04%  "
$ws.Range("D27").Value = "16.29"
$ws.Range("E27").Value = "  +3.75%  "
$ws.Range("E28").Value = "  +3.39%  "
$ws.Range("D29").Value = "6.87"
$ws.Range("E29").Value = "  +2.36%  "
$ws.Range("D30").Value = "0.983"
$ws.Range("E30").Value = "  -1.29%  "
$ws.Range("D31").Value = "3.79"
$ws.Range("E31").Value = "  +11.69%  "
$ws.Range("D32").Value = "0.0510"
$ws.Range("E32").Value = "  +2.58%  "
$ws.Range("D33").Value = "1.16"
$ws.Range("E33").Value = "  +4.53%  "
$ws.Range("E34").Value = "  +6.60%  "
$ws.Range("D35").Value = "1.549.59"
$ws.Range("E35").Value = "  +8.57%  "
$ws.Range("D36").Value = "1.76"
$ws.Range("E36").Value = "  +4.31%  "
$ws.Range("B37").Value = "Aave"
$ws.Range("C37").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D37").Value = "84.44"
$ws.Range("E37").Value = "  +9.21%  "
$ws.Range("B38").Value = "TrustWalletToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D38").Value = "1.03"
$ws.Range("E38").Value = "  +0.99%  "
$ws.Range("D39").Value = "0.618"
$ws.Range("E39").Value = "  +7.75%  "
$ws.Range("D40").Value = "0.0183"
$ws.Range("E40").Value = "  +5.12%  "
$ws.Range("D41").Value = "2.67"
$ws.Range("E41").Value = "  -0.87%  "
$ws.Range("D42").Value = "2.29"
$ws.Range("E42").Value = "  -0.40%  "
$ws.Range("D43").Value = "2.11"
$ws.Range("E43").Value = "  +7.92%  "
$ws.Range("E44").Value = "  +2.96%  "
$ws.Range("D45").Value = "0.0500"
$ws.Range("E45").Value = "  -0.30%  "
$ws.Range("D46").Value = "54.56"
$ws.Range("E46").Value = "  +7.45%  "
$ws.Range("E47").Value = "  +3.91%  "
$ws.Range("B48").Value = "RocketPoolETH"
$ws.Range("C48").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D48").Value = "1.877.31"
$ws.Range("E48").Value = "  +5.04%  "
$ws.Range("B49").Value = "PaxDollar"
$ws.Range("C49").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D49").Value = "0.986"
$ws.Range("E49").Value = "  -1.09%  "
$ws.Range("D50").Value = "5.58"
$ws.Range("E50").Value = "  +4.23%  "
$ws.Range("D51").Value = "94.34"
$ws.Range("E51").Value = "  +0.32%  "
